$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet to match updated export timestamp
$ws.Name = "牌局記錄0808_172450"

# Append 6 new game rounds (rows 869-874)
# Row 869
$ws.Cells.Item(869, 1).Value = 8
$ws.Cells.Item(869, 2).Value = "1,7"
$ws.Cells.Item(869, 3).Value = 7
$ws.Cells.Item(869, 4).Value = "7,K"
$ws.Cells.Item(869, 5).Value = 1
$ws.Cells.Item(869, 6).Value = "閒"
$ws.Cells.Item(869, 7).Value = "藍"
$ws.Cells.Item(869, 8).Value = "藍"
$ws.Cells.Item(869, 9).Value = "藍"
$ws.Cells.Item(869, 10).Value = "紅"
$ws.Cells.Item(869, 11).Value = "紅"
$ws.Cells.Item(869, 12).Value = "紅"

# Row 870
$ws.Cells.Item(870, 1).Value = 9
$ws.Cells.Item(870, 2).Value = "8,1"
$ws.Cells.Item(870, 3).Value = 2
$ws.Cells.Item(870, 4).Value = "4,8"
$ws.Cells.Item(870, 5).Value = 7
$ws.Cells.Item(870, 6).Value = "閒"
$ws.Cells.Item(870, 7).Value = "藍"
$ws.Cells.Item(870, 8).Value = "藍"
$ws.Cells.Item(870, 9).Value = "藍"
$ws.Cells.Item(870, 10).Value = "紅"
$ws.Cells.Item(870, 11).Value = "紅"
$ws.Cells.Item(870, 12).Value = "紅"

# Row 871
$ws.Cells.Item(871, 1).Value = 6
$ws.Cells.Item(871, 2).Value = "6,K"
$ws.Cells.Item(871, 3).Value = 1
$ws.Cells.Item(871, 4).Value = "5,9,7"
$ws.Cells.Item(871, 5).Value = 5
$ws.Cells.Item(871, 6).Value = "閒"
$ws.Cells.Item(871, 7).Value = "藍"
$ws.Cells.Item(871, 8).Value = "藍"
$ws.Cells.Item(871, 9).Value = "藍"
$ws.Cells.Item(871, 10).Value = "紅"
$ws.Cells.Item(871, 11).Value = "紅"
$ws.Cells.Item(871, 12).Value = "紅"

# Row 872
$ws.Cells.Item(872, 1).Value = 2
$ws.Cells.Item(872, 2).Value = "9,3"
$ws.Cells.Item(872, 3).Value = 9
$ws.Cells.Item(872, 4).Value = "9,K"
$ws.Cells.Item(872, 5).Value = 7
$ws.Cells.Item(872, 6).Value = "莊"
$ws.Cells.Item(872, 7).Value = "藍"
$ws.Cells.Item(872, 8).Value = "藍"
$ws.Cells.Item(872, 9).Value = "藍"
$ws.Cells.Item(872, 10).Value = "紅"
$ws.Cells.Item(872, 11).Value = "紅"
$ws.Cells.Item(872, 12).Value = "紅"

# Row 873
$ws.Cells.Item(873, 1).Value = 6
$ws.Cells.Item(873, 2).Value = "4,2"
$ws.Cells.Item(873, 3).Value = 1
$ws.Cells.Item(873, 4).Value = "K,4,7"
$ws.Cells.Item(873, 5).Value = 5
$ws.Cells.Item(873, 6).Value = "閒"
$ws.Cells.Item(873, 7).Value = "紅"
$ws.Cells.Item(873, 8).Value = "紅"
$ws.Cells.Item(873, 9).Value = "藍"
$ws.Cells.Item(873, 10).Value = "藍"
$ws.Cells.Item(873, 11).Value = "藍"
$ws.Cells.Item(873, 12).Value = "紅"

# Row 874
$ws.Cells.Item(874, 1).Value = 6
$ws.Cells.Item(874, 2).Value = "3,2,1"
$ws.Cells.Item(874, 3).Value = 9
$ws.Cells.Item(874, 4).Value = "3,J,6"
$ws.Cells.Item(874, 5).Value = 3
$ws.Cells.Item(874, 6).Value = "莊"
$ws.Cells.Item(874, 7).Value = "紅"
$ws.Cells.Item(874, 8).Value = "藍"
$ws.Cells.Item(874, 9).Value = "藍"
$ws.Cells.Item(874, 10).Value = "藍"
$ws.Cells.Item(874, 11).Value = "紅"
$ws.Cells.Item(874, 12).Value = "紅"
